$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.258986
$ws.Range("H2").Value = 9.776958
$ws.Range("I2").Value = 0.0928452675546778
$ws.Range("J2").Value = 0.09284526755467781
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04996866666666667
$ws.Range("N2").Value = 0.149906
$ws.Range("O2").Value = 0.06831051926220302
$ws.Range("P2").Value = 0.06831051926220301
$ws.Range("Q2").Value = 0.1628471851053334
$ws.Range("R2").Value = 1.465624665948
$ws.Range("S2").Value = 0.006342308437698212
$ws.Range("T2").Value = 0.006342308437698211
$ws.Range("G3").Value = 3.258986
$ws.Range("H3").Value = 9.776958
$ws.Range("I3").Value = 0.0928452675546778
$ws.Range("J3").Value = 0.09284526755467781
$ws.Range("O3").Value = 0.8063873019518528
$ws.Range("P3").Value = 0.8063873019518527
$ws.Range("Q3").Value = 1.922367208533333
$ws.Range("R3").Value = 17.3013048768
$ws.Range("S3").Value = 0.07486924480241453
$ws.Range("T3").Value = 0.07486924480241453
$ws.Range("G4").Value = 3.258986
$ws.Range("H4").Value = 9.776958
$ws.Range("I4").Value = 0.0928452675546778
$ws.Range("J4").Value = 0.09284526755467781
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.09165766666666668
$ws.Range("N4").Value = 0.274973
$ws.Range("O4").Value = 0.1253021787859442
$ws.Range("P4").Value = 0.1253021787859442
$ws.Range("Q4").Value = 0.2987110524593334
$ws.Range("R4").Value = 2.688399472134
$ws.Range("S4").Value = 0.01163371431456506
$ws.Range("T4").Value = 0.01163371431456506
$ws.Range("I5").Value = 0.5926110879358332
$ws.Range("J5").Value = 0.5926110879358333
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04996866666666667
$ws.Range("N5").Value = 0.149906
$ws.Range("O5").Value = 0.06831051926220302
$ws.Range("P5").Value = 0.06831051926220301
$ws.Range("Q5").Value = 1.039418056237778
$ws.Range("R5").Value = 9.35476250614
$ws.Range("S5").Value = 0.04048157113743582
$ws.Range("T5").Value = 0.04048157113743582
$ws.Range("I6").Value = 0.5926110879358332
$ws.Range("J6").Value = 0.5926110879358333
$ws.Range("O6").Value = 0.8063873019518528
$ws.Range("P6").Value = 0.8063873019518527
$ws.Range("S6").Value = 0.4778740563073287
$ws.Range("T6").Value = 0.4778740563073287
$ws.Range("I7").Value = 0.5926110879358332
$ws.Range("J7").Value = 0.5926110879358333
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.09165766666666668
$ws.Range("N7").Value = 0.274973
$ws.Range("O7").Value = 0.1253021787859442
$ws.Range("P7").Value = 0.1253021787859442
$ws.Range("Q7").Value = 1.906607481874445
$ws.Range("R7").Value = 17.15946733687
$ws.Range("S7").Value = 0.07425546049106868
$ws.Range("T7").Value = 0.07425546049106867
$ws.Range("G8").Value = 10.70406433333333
$ws.Range("H8").Value = 32.112193
$ws.Range("I8").Value = 0.3049481393755043
$ws.Range("J8").Value = 0.3049481393755044
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04996866666666667
$ws.Range("N8").Value = 0.149906
$ws.Range("O8").Value = 0.06831051926220302
$ws.Range("P8").Value = 0.06831051926220301
$ws.Range("Q8").Value = 0.5348678226508889
$ws.Range("R8").Value = 4.813810403858001
$ws.Range("S8").Value = 0.02083116574878336
$ws.Range("T8").Value = 0.02083116574878336
$ws.Range("G9").Value = 10.70406433333333
$ws.Range("H9").Value = 32.112193
$ws.Range("I9").Value = 0.3049481393755043
$ws.Range("J9").Value = 0.3049481393755044
$ws.Range("O9").Value = 0.8063873019518528
$ws.Range("P9").Value = 0.8063873019518527
$ws.Range("Q9").Value = 6.31397074808889
$ws.Range("R9").Value = 56.82573673280001
$ws.Range("S9").Value = 0.2459063073462505
$ws.Range("T9").Value = 0.2459063073462505
$ws.Range("G10").Value = 10.70406433333333
$ws.Range("H10").Value = 32.112193
$ws.Range("I10").Value = 0.3049481393755043
$ws.Range("J10").Value = 0.3049481393755044
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.09165766666666668
$ws.Range("N10").Value = 0.274973
$ws.Range("O10").Value = 0.1253021787859442
$ws.Range("P10").Value = 0.1253021787859442
$ws.Range("Q10").Value = 0.9811095606432225
$ws.Range("R10").Value = 8.829986045789003
$ws.Range("S10").Value = 0.03821066628047048
$ws.Range("T10").Value = 0.03821066628047047
$ws.Range("G11").Value = 0.3368143333333333
$ws.Range("H11").Value = 1.010443
$ws.Range("I11").Value = 0.009595505133984546
$ws.Range("J11").Value = 0.00959550513398455
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04996866666666667
$ws.Range("N11").Value = 0.149906
$ws.Range("O11").Value = 0.06831051926220302
$ws.Range("P11").Value = 0.06831051926220301
$ws.Range("Q11").Value = 0.01683016315088889
$ws.Range("R11").Value = 0.151471468358
$ws.Range("S11").Value = 0.0006554739382856194
$ws.Range("T11").Value = 0.0006554739382856195
$ws.Range("G12").Value = 0.3368143333333333
$ws.Range("H12").Value = 1.010443
$ws.Range("I12").Value = 0.009595505133984546
$ws.Range("J12").Value = 0.00959550513398455
$ws.Range("O12").Value = 0.8063873019518528
$ws.Range("P12").Value = 0.8063873019518527
$ws.Range("Q12").Value = 0.1986755480888889
$ws.Range("R12").Value = 1.7880799328
$ws.Range("S12").Value = 0.007737693495858951
$ws.Range("T12").Value = 0.007737693495858952
$ws.Range("G13").Value = 0.3368143333333333
$ws.Range("H13").Value = 1.010443
$ws.Range("I13").Value = 0.009595505133984546
$ws.Range("J13").Value = 0.00959550513398455
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09165766666666668
$ws.Range("N13").Value = 0.274973
$ws.Range("O13").Value = 0.1253021787859442
$ws.Range("P13").Value = 0.1253021787859442
$ws.Range("Q13").Value = 0.03087161589322223
$ws.Range("R13").Value = 0.277844543039
$ws.Range("S13").Value = 0.001202337699839977
$ws.Range("T13").Value = 0.001202337699839977
